$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove two providers from the list ("ACCESORIOS Y HERRAJES JM SAS" and
# "AGENCIA DE ADUANAS SERVICIOS"), keeping the remaining providers packed
# at the top of the table (rows 2-4) and leaving everything below the
# table (row 7 onward) untouched.

# Drop every existing hyperlink; we'll recreate only the two that survive.
$ws.Hyperlinks.Delete()

# New, compacted provider list (A2:B4) - values shifted up so the gap left
# by the two removed providers disappears.
$ws.Range("A2").Value = "ACDC ELECTRIC SAS"
$ws.Range("B2").Value = "japsequiposelectricos@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:japsequiposelectricos@gmail.com")
$ws.Range("B2").Style = "Hipervínculo"

$ws.Range("A3").Value = "PROVEEDOR SIN CORREO"
$ws.Range("B3").Clear()

$ws.Range("A4").Value = "PROVEEDOR SIN DOCUMENTO"
$ws.Range("B4").Value = "japsequiposelectricos@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:japsequiposelectricos@gmail.com")
$ws.Range("B4").Style = "Hipervínculo"

# Rows 5 and 6 no longer hold any provider - clear them out entirely so the
# now-unused rows disappear from the saved sheet.
$ws.Range("A5:B5").Clear()
$ws.Range("A6:B6").Clear()

# Reflect the resulting selection as seen after the edit.
$ws.Range("A2:B4").Select()
